# Apply the cryptos-list refresh described in the commit:
#   "Updated cryptos list on Fri May 19 16:19:45 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing text storage so that numeric-looking
# strings (e.g. "1.0000", "0.9997") are not reinterpreted as numbers by Excel.
function Set-TextValue {
    param($CellRef, $Text)
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.ClearFormats()
}

$ws.Range("D2").Value = '26.931.92'
$ws.Range("E2").Value = '  -0.64%  '
$ws.Range("D3").Value = '1.816.91'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("E4").Value = '  -0.13%  '
Set-TextValue "D5" '309.89'
$ws.Range("E5").Value = '  -0.63%  '
Set-TextValue "D6" '1.0000'
$ws.Range("E6").Value = '  -0.17%  '
Set-TextValue "D7" '0.4671'
$ws.Range("E7").Value = '  +1.39%  '
Set-TextValue "D8" '0.3698'
$ws.Range("E8").Value = '  -1.19%  '
Set-TextValue "D9" '0.07368'
$ws.Range("E9").Value = '  -0.41%  '
Set-TextValue "D10" '0.8708'
$ws.Range("E10").Value = '  +1.33%  '
Set-TextValue "D11" '20.42'
$ws.Range("E11").Value = '  -0.82%  '
$ws.Range("D12").Value = '1.840.23'
$ws.Range("E12").Value = '  +1.50%  '
Set-TextValue "D13" '5.365'
$ws.Range("E13").Value = '  -0.18%  '
Set-TextValue "D14" '92.25'
$ws.Range("E14").Value = '  +0.66%  '
$ws.Range("B15").Value = 'TRON'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue "D15" '0.07076'
$ws.Range("E15").Value = '  -0.22%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D16" '6.507'
$ws.Range("E16").Value = '  -2.17%  '
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("E20").Value = '  -0.91%  '
$ws.Range("D21").Value = '26.968.41'
$ws.Range("E21").Value = '  -0.59%  '
Set-TextValue "D22" '5.355'
$ws.Range("E22").Value = '  +0.66%  '
Set-TextValue "D23" '10.55'
$ws.Range("E23").Value = '  -2.91%  '
$ws.Range("D24").Value = '2.101.19'
$ws.Range("E24").Value = '  +2.61%  '
Set-TextValue "D25" '1.893'
$ws.Range("E25").Value = '  -1.78%  '
Set-TextValue "D26" '151.84'
Set-TextValue "D27" '2.201'
$ws.Range("E27").Value = '  -0.05%  '
Set-TextValue "D28" '18.37'
$ws.Range("E28").Value = '  -0.44%  '
Set-TextValue "D29" '5.309'
$ws.Range("E29").Value = '  +0.86%  '
Set-TextValue "D30" '115.60'
$ws.Range("E30").Value = '  -0.95%  '
$ws.Range("E31").Value = '  +0.19%  '
Set-TextValue "D32" '0.7656'
$ws.Range("E32").Value = '  -0.76%  '
$ws.Range("E33").Value = '  -0.42%  '
$ws.Range("E34").Value = '  -0.83%  '
Set-TextValue "D35" '2.920'
$ws.Range("E35").Value = '  +1.14%  '
Set-TextValue "D36" '0.9997'
$ws.Range("E36").Value = '  -0.16%  '
Set-TextValue "D37" '1.098'
$ws.Range("E37").Value = '  -2.49%  '
Set-TextValue "D38" '0.01962'
$ws.Range("E38").Value = '  +0.14%  '
Set-TextValue "D39" '0.05266'
$ws.Range("E39").Value = '  +0.73%  '
Set-TextValue "D40" '0.5395'
$ws.Range("E40").Value = '  +2.03%  '
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D41" '2.940'
$ws.Range("E41").Value = '  +0.65%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D42" '7.272'
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D43" '2.367'
$ws.Range("E43").Value = '  -0.83%  '
Set-TextValue "D44" '0.1663'
$ws.Range("E44").Value = '  -0.92%  '
Set-TextValue "D45" '8.489'
$ws.Range("E45").Value = '  -1.14%  '
Set-TextValue "D46" '0.4970'
$ws.Range("E46").Value = '  -1.23%  '
Set-TextValue "D47" '10.39'
$ws.Range("E47").Value = '  +0.22%  '
Set-TextValue "D48" '1.678'
$ws.Range("E48").Value = '  +0.58%  '
Set-TextValue "D49" '0.9995'
$ws.Range("E49").Value = '  -0.18%  '
Set-TextValue "D50" '103.27'
Set-TextValue "D51" '0.06278'
$ws.Range("E51").Value = '  -0.66%  '
